# Reproduces the crypto-price/volume refresh described in the diff.
# D (Price) and E (Volume(1h)) columns are stored as plain text in the
# workbook (inline strings), including values that look like plain
# decimals (e.g. "4.60"). Excel normally auto-converts a bare decimal
# typed into a cell into a Number (and would drop the trailing zero),
# so those specific values are written with a leading apostrophe to
# force a text entry, matching the source data exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '64.052.41'
$ws.Cells.Item(2, 5).Value = '  -2.71%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '2.629.86'
$ws.Cells.Item(3, 5).Value = '  -1.05%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.01%  '

# Row 5
$ws.Cells.Item(5, 4).Value = "'577.63"
$ws.Cells.Item(5, 5).Value = '  -3.50%  '

# Row 6
$ws.Cells.Item(6, 4).Value = "'157.22"
$ws.Cells.Item(6, 5).Value = '  -0.15%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  +0.06%  '

# Row 8
$ws.Cells.Item(8, 4).Value = "'0.631"
$ws.Cells.Item(8, 5).Value = '  +0.07%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  -4.82%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  +0.33%  '

# Row 11
$ws.Cells.Item(11, 4).Value = "'0.386"
$ws.Cells.Item(11, 5).Value = '  -3.01%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  -0.23%  '

# Row 13
$ws.Cells.Item(13, 4).Value = "'28.48"
$ws.Cells.Item(13, 5).Value = '  -0.69%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '3.106.13'
$ws.Cells.Item(14, 5).Value = '  -0.89%  '

# Row 15
$ws.Cells.Item(15, 5).Value = '  -5.99%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '63.862.94'
$ws.Cells.Item(16, 5).Value = '  -2.73%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '2.623.61'
$ws.Cells.Item(17, 5).Value = '  -0.78%  '

# Row 18
$ws.Cells.Item(18, 5).Value = '  -3.45%  '

# Row 19
$ws.Cells.Item(19, 4).Value = "'7.71"
$ws.Cells.Item(19, 5).Value = '  +3.18%  '

# Row 20
$ws.Cells.Item(20, 4).Value = "'4.60"
$ws.Cells.Item(20, 5).Value = '  -3.04%  '

# Row 21
$ws.Cells.Item(21, 4).Value = "'345.01"
$ws.Cells.Item(21, 5).Value = '  -1.44%  '

# Row 22
$ws.Cells.Item(22, 5).Value = '  -0.03%  '

# Row 23
$ws.Cells.Item(23, 4).Value = "'67.60"
$ws.Cells.Item(23, 5).Value = '  -2.32%  '

# Row 24
$ws.Cells.Item(24, 5).Value = '  +0.91%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  -1.89%  '

# Row 26
$ws.Cells.Item(26, 4).Value = "'600.51"
$ws.Cells.Item(26, 5).Value = '  +8.24%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  -4.06%  '

# Row 28
$ws.Cells.Item(28, 5).Value = '  -1.08%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  -0.77%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  +0.05%  '

# Row 31
$ws.Cells.Item(31, 4).Value = "'7.93"
$ws.Cells.Item(31, 5).Value = '  +0.13%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  -2.50%  '

# Row 33
$ws.Cells.Item(33, 4).Value = "'1.74"
$ws.Cells.Item(33, 5).Value = '  -1.78%  '

# Row 34
$ws.Cells.Item(34, 4).Value = "'6.63"
$ws.Cells.Item(34, 5).Value = '  +1.93%  '

# Row 35
$ws.Cells.Item(35, 4).Value = "'5.36"
$ws.Cells.Item(35, 5).Value = '  -1.27%  '

# Row 36
$ws.Cells.Item(36, 4).Value = "'0.410"

# Row 37
$ws.Cells.Item(37, 4).Value = "'19.89"
$ws.Cells.Item(37, 5).Value = '  -2.36%  '

# Row 38
$ws.Cells.Item(38, 4).Value = "'0.999"
$ws.Cells.Item(38, 5).Value = '  -0.03%  '

# Row 39
$ws.Cells.Item(39, 4).Value = "'154.63"
$ws.Cells.Item(39, 5).Value = '  -0.21%  '

# Row 40
$ws.Cells.Item(40, 5).Value = '  -2.54%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  -0.03%  '

# Row 42
$ws.Cells.Item(42, 4).Value = "'41.57"
$ws.Cells.Item(42, 5).Value = '  -2.14%  '

# Row 43
$ws.Cells.Item(43, 5).Value = '  +5.87%  '

# Row 44
$ws.Cells.Item(44, 4).Value = "'157.20"
$ws.Cells.Item(44, 5).Value = '  -2.98%  '

# Row 45
$ws.Cells.Item(45, 5).Value = '  -3.00%  '

# Row 46
$ws.Cells.Item(46, 4).Value = "'23.28"
$ws.Cells.Item(46, 5).Value = '  +2.89%  '

# Row 47
$ws.Cells.Item(47, 5).Value = '  -1.08%  '

# Row 48
$ws.Cells.Item(48, 2).Value = 'Mantle'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(48, 4).Value = "'0.631"
$ws.Cells.Item(48, 5).Value = '  -1.06%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'Stellar'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(49, 4).Value = "'0.102"
$ws.Cells.Item(49, 5).Value = '  +1.79%  '

# Row 50
$ws.Cells.Item(50, 5).Value = '  -2.14%  '

# Row 51
$ws.Cells.Item(51, 4).Value = "'19.06"
$ws.Cells.Item(51, 5).Value = '  -3.68%  '
